# Applies the "Add files via upload" edit to the "peak" sheet:
#   - Rename the two shared-string header labels used by columns Q/R
#     from the "...SmA03S02..." family to "...RaA03S02..." (A03S02 LPV/LPD).
#   - Update row 2 (the single data row) with refreshed peak-date figures:
#     a handful of previously-blank LPV-style date cells now carry values,
#     R2's date moves, and AE2's count drops by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("peak")

# Header text fix-ups (row 1)
$ws.Range("Q1").Value = "DayDeaMeRaA03S02_LPV"
$ws.Range("R1").Value = "DayDeaMeRaA03S02_LPD"

# Row 2 data fix-ups
$ws.Range("R2").Value = 44448
$ws.Range("T2").Value = 44444
$ws.Range("V2").Value = 44441
$ws.Range("X2").Value = 44442
$ws.Range("Z2").Value = 44438
$ws.Range("AB2").Value = 44434
$ws.Range("AD2").Value = 44431
$ws.Range("AE2").Value = 33235219
$ws.Range("AF2").Value = 44432
